$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the table with a new "2022" column (S), mirroring the formatting
# of the existing "2021" column (R) for the header, share and volume rows.
$ws.Range("R4:R6").Copy()
$ws.Range("S4:S6").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("S4").Value = 2022
$ws.Range("S5").Value = 30
$ws.Range("S6").Value = 11928.6

# Update the active selection to match the edited workbook state.
$ws.Range("T3").Select()
